# Apply updated crypto price/volume data to cryptos.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($cellRef, $value) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $value
}

# Row 2
Set-CellText 'D2' '68.750.28'
Set-CellText 'E2' '  -0.63%  '

# Row 3
Set-CellText 'D3' '3.836.77'
Set-CellText 'E3' '  +2.29%  '

# Row 4
Set-CellText 'E4' '  -0.01%  '

# Row 5
Set-CellText 'D5' '601.07'
Set-CellText 'E5' '  -0.23%  '

# Row 6
Set-CellText 'D6' '161.53'
Set-CellText 'E6' '  -3.26%  '

# Row 7
Set-CellText 'D7' '3.833.37'
Set-CellText 'E7' '  +2.23%  '

# Row 8
Set-CellText 'E8' '  +0.02%  '

# Row 9
Set-CellText 'E9' '  -1.55%  '

# Row 10
Set-CellText 'E10' '  -1.14%  '

# Row 11
Set-CellText 'D11' '6.29'
Set-CellText 'E11' '  -1.41%  '

# Row 12
Set-CellText 'E12' '  -0.24%  '

# Row 13
Set-CellText 'D13' '36.83'
Set-CellText 'E13' '  -3.20%  '

# Row 14
Set-CellText 'E14' '  -2.27%  '

# Row 15
Set-CellText 'D15' '4.481.77'
Set-CellText 'E15' '  +2.29%  '

# Row 16
Set-CellText 'D16' '3.819.04'
Set-CellText 'E16' '  +1.86%  '

# Row 17
Set-CellText 'D17' '68.867.06'
Set-CellText 'E17' '  -0.48%  '

# Row 18
Set-CellText 'D18' '7.51'
Set-CellText 'E18' '  +1.65%  '

# Row 20
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-CellText 'D20' '17.10'
Set-CellText 'E20' '  -2.05%  '

# Row 21
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-CellText 'D21' '11.30'
Set-CellText 'E21' '  +1.28%  '

# Row 22
Set-CellText 'D22' '483.56'
Set-CellText 'E22' '  -2.22%  '

# Row 23
Set-CellText 'D23' '0.718'
Set-CellText 'E23' '  -1.34%  '

# Row 24
Set-CellText 'E24' '  +2.88%  '

# Row 25
Set-CellText 'D25' '83.94'
Set-CellText 'E25' '  -1.11%  '

# Row 26
Set-CellText 'E26' '  -2.66%  '

# Row 27
Set-CellText 'D27' '12.07'
Set-CellText 'E27' '  -2.00%  '

# Row 28
Set-CellText 'E28' '  -0.12%  '

# Row 29
Set-CellText 'D29' '9.98'
Set-CellText 'E29' '  -1.31%  '

# Row 30
Set-CellText 'E30' '  -1.31%  '

# Row 31
Set-CellText 'D31' '7.92'
Set-CellText 'E31' '  -1.95%  '

# Row 32
Set-CellText 'D32' '3.987.98'
Set-CellText 'E32' '  +2.33%  '

# Row 33
Set-CellText 'E33' '  -4.02%  '

# Row 34
Set-CellText 'D34' '32.07'
Set-CellText 'E34' '  +1.41%  '

# Row 35
Set-CellText 'D35' '3.784.06'
Set-CellText 'E35' '  +2.54%  '

# Row 36
Set-CellText 'E36' '  -1.78%  '

# Row 37
Set-CellText 'E37' '  +0.87%  '

# Row 38
Set-CellText 'E38' '  +3.23%  '

# Row 39
Set-CellText 'E39' '  -1.55%  '

# Row 40
Set-CellText 'D40' '0.999'
Set-CellText 'E40' '  -0.06%  '

# Row 41
Set-CellText 'E41' '  -1.99%  '

# Row 42
Set-CellText 'D42' '436.35'
Set-CellText 'E42' '  +1.48%  '

# Row 43
Set-CellText 'E43' '  -1.01%  '

# Row 44
Set-CellText 'E44' '  -0.68%  '

# Row 45
Set-CellText 'E45' '  -0.62%  '

# Row 46
Set-CellText 'E46' '  -0.03%  '

# Row 47
Set-CellText 'E47' '  -1.23%  '

# Row 48
Set-CellText 'D48' '143.10'
Set-CellText 'E48' '  +1.30%  '

# Row 49
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-CellText 'D49' '2.822.19'
Set-CellText 'E49' '  +0.91%  '

# Row 50
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-CellText 'D50' '0.0360'
Set-CellText 'E50' '  +2.06%  '

# Row 51
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-CellText 'D51' '26.10'
Set-CellText 'E51' '  +11.80%  '

Write-Output "Applied cryptos.xlsx update."